$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,("G2", 0.1427386666666667)
    ,("H2", 0.428216)
    ,("I2", 0.001645492988648044)
    ,("J2", 0.001657527388465106)
    ,("M2", 9.363528666666667)
    ,("N2", 28.090586)
    ,("O2", 0.04175743631338733)
    ,("P2", 0.04324026421082073)
    ,("Q2", 1.336537597175111)
    ,("R2", 12.028838374576)
    ,("S2", 0.00006871156867759609)
    ,("T2", 0.00007167192221390288)
    ,("G3", 0.1427386666666667)
    ,("H3", 0.428216)
    ,("I3", 0.001645492988648044)
    ,("J3", 0.001657527388465106)
    ,("O3", 0.1749266505387075)
    ,("P3", 0.1811383852696593)
    ,("Q3", 5.598908022951111)
    ,("R3", 50.39017220656)
    ,("S3", 0.0002878405769891297)
    ,("T3", 0.0003002418346868046)
    ,("G4", 0.1427386666666667)
    ,("H4", 0.428216)
    ,("I4", 0.001645492988648044)
    ,("J4", 0.001657527388465106)
    ,("M4", 71.284935)
    ,("N4", 213.854805)
    ,("O4", 0.3179011075133629)
    ,("P4", 0.3291899382573772)
    ,("Q4", 10.17511657532)
    ,("R4", 91.57604917787999)
    ,("S4", 0.0005231040434966868)
    ,("T4", 0.00054564133866874)
    ,("G5", 0.1427386666666667)
    ,("H5", 0.428216)
    ,("I5", 0.001645492988648044)
    ,("J5", 0.001657527388465106)
    ,("M5", 23.0690325)
    ,("N5", 46.138065)
    ,("O5", 0.1028782726814826)
    ,("P5", 0.07102102180339065)
    ,("Q5", 3.292842940339999)
    ,("R5", 19.75705764204)
    ,("S5", 0.0001692854763816012)
    ,("T5", 0.0001177192887958975)
    ,("G6", 0.1427386666666667)
    ,("H6", 0.428216)
    ,("I6", 0.001645492988648044)
    ,("J6", 0.001657527388465106)
    ,("M6", 81.293813)
    ,("N6", 243.881439)
    ,("O6", 0.3625365329530597)
    ,("P6", 0.3754103904587522)
    ,("Q6", 11.60377047586933)
    ,("R6", 104.433934282824)
    ,("S6", 0.0005965513231030304)
    ,("T6", 0.0006222530040997613)
    ,("H7", 38.366549)
    ,("I7", 0.1474300058337887)
    ,("J7", 0.1485082429624034)
    ,("M7", 9.363528666666667)
    ,("N7", 28.090586)
    ,("O7", 0.04175743631338733)
    ,("P7", 0.04324026421082073)
    ,("Q7", 119.7487604675238)
    ,("R7", 1077.738844207714)
    ,("S7", 0.006156299079286752)
    ,("T7", 0.006421535663179081)
    ,("H8", 38.366549)
    ,("I8", 0.1474300058337887)
    ,("J8", 0.1485082429624034)
    ,("O8", 0.1749266505387075)
    ,("P8", 0.1811383852696593)
    ,("Q8", 501.6411787720378)
    ,("R8", 4514.77060894834)
    ,("S8", 0.02578943710940675)
    ,("T8", 0.02690054332944399)
    ,("H9", 38.366549)
    ,("I9", 0.1474300058337887)
    ,("J9", 0.1485082429624034)
    ,("M9", 71.284935)
    ,("N9", 213.854805)
    ,("O9", 0.3179011075133629)
    ,("P9", 0.3291899382573772)
    ,("Q9", 911.652317213105)
    ,("R9", 8204.870854917945)
    ,("S9", 0.04686816213526297)
    ,("T9", 0.04888741933150514)
    ,("H10", 38.366549)
    ,("I10", 0.1474300058337887)
    ,("J10", 0.1485082429624034)
    ,("M10", 23.0690325)
    ,("N10", 46.138065)
    ,("O10", 0.1028782726814826)
    ,("P10", 0.07102102180339065)
    ,("Q10", 295.0263885979475)
    ,("R10", 1770.158331587685)
    ,("S10", 0.01516734434160108)
    ,("T10", 0.01054720716141609)
    ,("H11", 38.366549)
    ,("I11", 0.1474300058337887)
    ,("J11", 0.1485082429624034)
    ,("M11", 81.293813)
    ,("N11", 243.881439)
    ,("O11", 0.3625365329530597)
    ,("P11", 0.3754103904587522)
    ,("Q11", 1039.654353287112)
    ,("R11", 9356.88917958401)
    ,("S11", 0.05344876316823111)
    ,("T11", 0.05575153747685909)
    ,("G12", 33.975493)
    ,("H12", 101.926479)
    ,("I12", 0.3916698735032837)
    ,("J12", 0.3945343717944063)
    ,("M12", 9.363528666666667)
    ,("N12", 28.090586)
    ,("O12", 0.04175743631338733)
    ,("P12", 0.04324026421082073)
    ,("Q12", 318.1305026696327)
    ,("R12", 2863.174524026694)
    ,("S12", 0.01635512979868584)
    ,("T12", 0.01705977047664031)
    ,("G13", 33.975493)
    ,("H13", 101.926479)
    ,("I13", 0.3916698735032837)
    ,("J13", 0.3945343717944063)
    ,("O13", 0.1749266505387075)
    ,("P13", 0.1811383852696593)
    ,("Q13", 1332.684862369127)
    ,("R13", 11994.16376132214)
    ,("S13", 0.06851349908884866)
    ,("T13", 0.07146531904021816)
    ,("G14", 33.975493)
    ,("H14", 101.926479)
    ,("I14", 0.3916698735032837)
    ,("J14", 0.3945343717944063)
    ,("M14", 71.284935)
    ,("N14", 213.854805)
    ,("O14", 0.3179011075133629)
    ,("P14", 0.3291899382573772)
    ,("Q14", 2421.940810097955)
    ,("R14", 21797.4672908816)
    ,("S14", 0.1245122865663126)
    ,("T14", 0.1298767454914137)
    ,("G15", 33.975493)
    ,("H15", 101.926479)
    ,("I15", 0.3916698735032837)
    ,("J15", 0.3945343717944063)
    ,("M15", 23.0690325)
    ,("N15", 46.138065)
    ,("O15", 0.1028782726814826)
    ,("P15", 0.07102102180339065)
    ,("Q15", 783.7817522205224)
    ,("R15", 4702.690513323135)
    ,("S15", 0.04029432004739262)
    ,("T15", 0.02802023422139756)
    ,("G16", 33.975493)
    ,("H16", 101.926479)
    ,("I16", 0.3916698735032837)
    ,("J16", 0.3945343717944063)
    ,("M16", 81.293813)
    ,("N16", 243.881439)
    ,("O16", 0.3625365329530597)
    ,("P16", 0.3754103904587522)
    ,("Q16", 2761.997374524809)
    ,("R16", 24857.97637072328)
    ,("S16", 0.1419946380020439)
    ,("T16", 0.1481123025647366)
    ,("G17", 1.889429)
    ,("H17", 3.778858)
    ,("I17", 0.02178135921157747)
    ,("J17", 0.01462710555448763)
    ,("M17", 9.363528666666667)
    ,("N17", 28.090586)
    ,("O17", 0.04175743631338733)
    ,("P17", 0.04324026421082073)
    ,("Q17", 17.69172260513133)
    ,("R17", 106.150335630788)
    ,("S17", 0.0009095337200964586)
    ,("T17", 0.0006324799088156084)
    ,("G18", 1.889429)
    ,("H18", 3.778858)
    ,("I18", 0.02178135921157747)
    ,("J18", 0.01462710555448763)
    ,("O18", 0.1749266505387075)
    ,("P18", 0.1811383852696593)
    ,("Q18", 74.11263838971333)
    ,("R18", 444.67583033828)
    ,("S18", 0.003810140211061668)
    ,("T18", 0.002649530281308753)
    ,("G19", 1.889429)
    ,("H19", 3.778858)
    ,("I19", 0.02178135921157747)
    ,("J19", 0.01462710555448763)
    ,("M19", 71.284935)
    ,("N19", 213.854805)
    ,("O19", 0.3179011075133629)
    ,("P19", 0.3291899382573772)
    ,("Q19", 134.687823452115)
    ,("R19", 808.1269407126899)
    ,("S19", 0.006924318216506867)
    ,("T19", 0.004815095974365921)
    ,("G20", 1.889429)
    ,("H20", 3.778858)
    ,("I20", 0.02178135921157747)
    ,("J20", 0.01462710555448763)
    ,("M20", 23.0690325)
    ,("N20", 46.138065)
    ,("O20", 0.1028782726814826)
    ,("P20", 0.07102102180339065)
    ,("Q20", 43.5872990074425)
    ,("R20", 174.34919602977)
    ,("S20", 0.002240828612341989)
    ,("T20", 0.001038831982505762)
    ,("G21", 1.889429)
    ,("H21", 3.778858)
    ,("I21", 0.02178135921157747)
    ,("J21", 0.01462710555448763)
    ,("M21", 81.293813)
    ,("N21", 243.881439)
    ,("O21", 0.3625365329530597)
    ,("P21", 0.3754103904587522)
    ,("Q21", 153.598887802777)
    ,("R21", 921.5933268166619)
    ,("S21", 0.007896538451570486)
    ,("T21", 0.005491167407491583)
    ,("G22", 37.948719)
    ,("H22", 113.846157)
    ,("I22", 0.4374732684627022)
    ,("J22", 0.4406727523002374)
    ,("M22", 9.363528666666667)
    ,("N22", 28.090586)
    ,("O22", 0.04175743631338733)
    ,("P22", 0.04324026421082073)
    ,("Q22", 355.3339182197781)
    ,("R22", 3198.005263978002)
    ,("S22", 0.01826776214664069)
    ,("T22", 0.01905480623997182)
    ,("G23", 37.948719)
    ,("H23", 113.846157)
    ,("I23", 0.4374732684627022)
    ,("J23", 0.4406727523002374)
    ,("O23", 0.1749266505387075)
    ,("P23", 0.1811383852696593)
    ,("Q23", 1488.53420192018)
    ,("R23", 13396.80781728162)
    ,("S23", 0.07652573355240126)
    ,("T23", 0.07982275078400154)
    ,("G24", 37.948719)
    ,("H24", 113.846157)
    ,("I24", 0.4374732684627022)
    ,("J24", 0.4406727523002374)
    ,("M24", 71.284935)
    ,("N24", 213.854805)
    ,("O24", 0.3179011075133629)
    ,("P24", 0.3291899382573772)
    ,("Q24", 2705.171967248265)
    ,("R24", 24346.54770523438)
    ,("S24", 0.1390732365517838)
    ,("T24", 0.1450650361214236)
    ,("G25", 37.948719)
    ,("H25", 113.846157)
    ,("I25", 0.4374732684627022)
    ,("J25", 0.4406727523002374)
    ,("M25", 23.0690325)
    ,("N25", 46.138065)
    ,("O25", 0.1028782726814826)
    ,("P25", 0.07102102180339065)
    ,("Q25", 875.4402319443675)
    ,("R25", 5252.641391666205)
    ,("S25", 0.04500649420376533)
    ,("T25", 0.03129702914927533)
    ,("G26", 37.948719)
    ,("H26", 113.846157)
    ,("I26", 0.4374732684627022)
    ,("J26", 0.4406727523002374)
    ,("M26", 81.293813)
    ,("N26", 243.881439)
    ,("O26", 0.3625365329530597)
    ,("P26", 0.3754103904587522)
    ,("Q26", 3084.996065975547)
    ,("R26", 27764.96459377993)
    ,("S26", 0.1586000420081112)
    ,("T26", 0.1654331300055651)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

"Applied $($updates.Count) cell updates"